$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(590).Insert()

$ws.Cells.Item(590, 1).Value = 5
$ws.Cells.Item(590, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(590, 3).Value = "Maule"
$ws.Cells.Item(590, 4).Value = 45275
$ws.Cells.Item(590, 5).Value = 7
$ws.Cells.Item(590, 6).Value = 100114014
$ws.Cells.Item(590, 7).Value = "Betarraga"
$ws.Cells.Item(590, 8).Value = "Sin especificar"
$ws.Cells.Item(590, 9).Value = "Primera"
$ws.Cells.Item(590, 10).Value = 5000
$ws.Cells.Item(590, 11).Value = 700
$ws.Cells.Item(590, 12).Value = 700
$ws.Cells.Item(590, 13).Value = 700
$ws.Cells.Item(590, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(590, 15).Value = "Región del Maule"
$ws.Cells.Item(590, 16).Value = 140
$ws.Cells.Item(590, 17).Value = 5
$ws.Cells.Item(590, 18).Value = "Hortaliza"
